# Insert a new weekly price record for "Coliflor" (Vega Modelo de Temuco)
# at row 665, pushing all existing rows 665-697 down to 666-698.
# This mirrors Excel's native "insert row" behaviour (Rows.Item(n).Insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 665; everything below shifts down by one.
$ws.Rows.Item(665).Insert()

# Populate the newly inserted row 665 with the new record's data.
$ws.Cells.Item(665, 1).Value  = 10
$ws.Cells.Item(665, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(665, 3).Value  = "La Araucanía"
$ws.Cells.Item(665, 4).Value  = 45147
$ws.Cells.Item(665, 5).Value  = 9
$ws.Cells.Item(665, 6).Value  = 100112008
$ws.Cells.Item(665, 7).Value  = "Coliflor"
$ws.Cells.Item(665, 8).Value  = "Sin especificar"
$ws.Cells.Item(665, 9).Value  = "Primera"
$ws.Cells.Item(665, 10).Value = 650
$ws.Cells.Item(665, 11).Value = 1000
$ws.Cells.Item(665, 12).Value = 1000
$ws.Cells.Item(665, 13).Value = 1000
$ws.Cells.Item(665, 14).Value = "$/unidad"
$ws.Cells.Item(665, 15).Value = "Región del Maule"
$ws.Cells.Item(665, 16).Value = 1000
$ws.Cells.Item(665, 17).Value = 1
$ws.Cells.Item(665, 18).Value = "Hortaliza"

# Keep the date column's display format consistent with the rest of column D.
$ws.Cells.Item(665, 4).NumberFormat = $ws.Cells.Item(666, 4).NumberFormat
